$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marker rename: TDY1948 -> TDY1970 (rows 41-43, column E)
$ws.Range("E41").Value = "TDY1970"
$ws.Range("E42").Value = "TDY1970"
$ws.Range("E43").Value = "TDY1970"

# marker_1 threshold reclassification: G418 -> NAT (rows 41-43, column J)
$ws.Range("J41").Value = "NAT"
$ws.Range("J42").Value = "NAT"
$ws.Range("J43").Value = "NAT"

# Update the active selection to reflect the last-edited cell
$ws.Range("E43").Select()
